$wb = $excel.ActiveWorkbook
$ds = $wb.Worksheets.Item("data")

# Add new "metadata" worksheet at the end of the workbook (after the last sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "metadata"

# Header row (row 1) - bold/bordered style like the "data" sheet headers
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Copy the header style (font/border/alignment) from the "data" sheet headers
$ds.Range("B1:F1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Data row (row 2)
$ws.Range("A2").Value = 0
$ds.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B2").Value = "Tuberous Sclerosis_Focal Cortical Dysplasia_Hemimegalencephaly"
$ws.Range("C2").Value = 20

# Force D2 ("0.41") to be stored as text rather than a number
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.41"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "2021-03-22T10:10:41.483099Z"
$ws.Range("F2").Value = "2021-10-05 14:35:52.134718"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/20/?format=json"

$excel.CutCopyMode = $false

# Update query-time timestamps on "data" sheet (column F, rows 2-12)
$ds.Range("F2").Value = "2021-10-05 14:35:52.138321"
$ds.Range("F3").Value = "2021-10-05 14:35:52.138329"
$ds.Range("F4").Value = "2021-10-05 14:35:52.138332"
$ds.Range("F5").Value = "2021-10-05 14:35:52.138335"
$ds.Range("F6").Value = "2021-10-05 14:35:52.138337"
$ds.Range("F7").Value = "2021-10-05 14:35:52.138340"
$ds.Range("F8").Value = "2021-10-05 14:35:52.138343"
$ds.Range("F9").Value = "2021-10-05 14:35:52.138345"
$ds.Range("F10").Value = "2021-10-05 14:35:52.138348"
$ds.Range("F11").Value = "2021-10-05 14:35:52.138351"
$ds.Range("F12").Value = "2021-10-05 14:35:52.138353"
